# SCHRÄG BOM update — smaller packages
# D3: zener diode package swap SOD-123 -> SOT-23-3, part swap BZT52-C4V7X -> PLVA650A,215
# RV6/RV7/RV8: trimmer pot swap 3296W-1-xxxLF -> 3266Y-1-xxxLF (smaller package)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- D3 (row 11): Diode, zener ---
$ws.Range("C11").Value = "SOT-23-3"
$ws.Range("H11").Value = "PLVA650A,215"
$ws.Range("I11").Value = "771-PLVA650A-T/R"

# --- RV6 (row 34): Multiturn trimmer pot, 100K ---
$ws.Range("H34").Value = "3266Y-1-104LF"
$ws.Range("I34").Value = "652-3266Y-1-104LF"

# --- RV7 (row 35): Multiturn trimmer pot, 5K ---
$ws.Range("H35").Value = "3266Y-1-502LF"
$ws.Range("I35").Value = "652-3266Y-1-502LF"

# --- RV8 (row 36): Multiturn trimmer pot, 250K ---
$ws.Range("H36").Value = "3266Y-1-254LF"
$ws.Range("I36").Value = "652-3266Y-1-254LF"

# --- view state: zoom out to 90%, reset selection to A1 ---
$excel.ActiveWindow.Zoom = 90
$null = $ws.Range("A1").Select()
